$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price value would otherwise be auto-coerced to a Number
# by Excel (losing the literal text formatting, e.g. trailing zeros).
# Force them to keep Text format before assigning the value.
$textCells = @("D5","D6","D8","D9","D11","D13","D15","D16","D20","D21","D22","D23","D26","D27","D28","D30","D31","D32","D36","D37","D41","D44","D45","D46","D47","D48","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = '37.028.08'
$ws.Range("E2").Value = '  +1.54%  '

# Row 3
$ws.Range("D3").Value = '2.054.00'
$ws.Range("E3").Value = '  -1.94%  '

# Row 4
$ws.Range("E4").Value = '  -0.01%  '

# Row 5
$ws.Range("D5").Value = '249.14'
$ws.Range("E5").Value = '  -0.71%  '

# Row 6
$ws.Range("D6").Value = '0.680'
$ws.Range("E6").Value = '  +3.96%  '

# Row 7
$ws.Range("E7").Value = '  -0.08%  '

# Row 8
$ws.Range("D8").Value = '54.26'
$ws.Range("E8").Value = '  +14.65%  '

# Row 9
$ws.Range("D9").Value = '60.59'
$ws.Range("E9").Value = '  +1.99%  '

# Row 10
$ws.Range("E10").Value = '  +1.80%  '

# Row 11
$ws.Range("D11").Value = '0.0787'
$ws.Range("E11").Value = '  +6.17%  '

# Row 12
$ws.Range("E12").Value = '  +6.22%  '

# Row 13
$ws.Range("D13").Value = '14.85'
$ws.Range("E13").Value = '  +3.05%  '

# Row 14
$ws.Range("D14").Value = '2.351.64'
$ws.Range("E14").Value = '  -1.97%  '

# Row 15
$ws.Range("D15").Value = '0.814'
$ws.Range("E15").Value = '  -1.13%  '

# Row 16
$ws.Range("D16").Value = '5.23'
$ws.Range("E16").Value = '  +3.18%  '

# Row 17
$ws.Range("E17").Value = '  -1.93%  '

# Row 18
$ws.Range("D18").Value = '36.973.51'
$ws.Range("E18").Value = '  +1.47%  '

# Row 19
$ws.Range("D19").Value = '0.0₃0928'
$ws.Range("E19").Value = '  +12.45%  '

# Row 20
$ws.Range("D20").Value = '72.64'
$ws.Range("E20").Value = '  +0.28%  '

# Row 21
$ws.Range("D21").Value = '14.17'
$ws.Range("E21").Value = '  +7.33%  '

# Row 22
$ws.Range("D22").Value = '5.34'
$ws.Range("E22").Value = '  +4.34%  '

# Row 23
$ws.Range("D23").Value = '235.80'
$ws.Range("E23").Value = '  -1.28%  '

# Row 25
$ws.Range("E25").Value = '  -2.02%  '

# Row 26
$ws.Range("D26").Value = '170.14'
$ws.Range("E26").Value = '  -0.01%  '

# Row 27
$ws.Range("D27").Value = '8.98'
$ws.Range("E27").Value = '  -0.91%  '

# Row 28
$ws.Range("D28").Value = '20.01'
$ws.Range("E28").Value = '  -5.49%  '

# Row 29
$ws.Range("E29").Value = '  +0.47%  '

# Row 30
$ws.Range("D30").Value = '0.125'
$ws.Range("E30").Value = '  +2.50%  '

# Row 31
$ws.Range("D31").Value = '4.55'
$ws.Range("E31").Value = '  +2.97%  '

# Row 32
$ws.Range("D32").Value = '0.0620'
$ws.Range("E32").Value = '  +1.49%  '

# Row 33
$ws.Range("E33").Value = '  +7.47%  '

# Row 34
$ws.Range("E34").Value = '  +6.73%  '

# Row 35
$ws.Range("E35").Value = '  +0.02%  '

# Row 36
$ws.Range("D36").Value = '0.0860'
$ws.Range("E36").Value = '  -7.57%  '

# Row 37
$ws.Range("D37").Value = '2.28'
$ws.Range("E37").Value = '  -1.84%  '

# Row 38
$ws.Range("E38").Value = '  -6.11%  '

# Row 39
$ws.Range("E39").Value = '  +0.30%  '

# Row 40
$ws.Range("E40").Value = '  +23.55%  '

# Row 41
$ws.Range("D41").Value = '17.78'
$ws.Range("E41").Value = '  +11.48%  '

# Row 42
$ws.Range("E42").Value = '  +0.62%  '

# Row 43
$ws.Range("E43").Value = '  -1.80%  '

# Row 44
$ws.Range("D44").Value = '96.08'
$ws.Range("E44").Value = '  -1.02%  '

# Row 45
$ws.Range("D45").Value = '2.79'
$ws.Range("E45").Value = '  +1.48%  '

# Row 46
$ws.Range("D46").Value = '4.12'
$ws.Range("E46").Value = '  +53.93%  '

# Row 47
$ws.Range("B47").Value = 'Gas'
$ws.Range("C47").Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range("D47").Value = '13.37'
$ws.Range("E47").Value = '  -52.84%  '

# Row 48
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").Value = '2.38'
$ws.Range("E48").Value = '  +7.55%  '

# Row 49
$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").Value = '1.289.85'
$ws.Range("E49").Value = '  -2.83%  '

# Row 50
$ws.Range("B50").Value = 'MXToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D50").Value = '2.92'
$ws.Range("E50").Value = '  +2.69%  '

# Row 51
$ws.Range("D51").Value = '4.07'
